# Updated cryptos list - refresh Price (D) and Volume(1h) (E) columns
# Values are written with a leading apostrophe to force text storage
# (matching the source data's text-typed cells), then the cell style is
# reset to Normal so no stray number-format style is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.017.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.66%  "

$ws.Range("D3").Value = "'2.300.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.05%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'252.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.55%  "

$ws.Range("D6").Value = "'0.643"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.89%  "

$ws.Range("D7").Value = "'74.12"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.15%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("E9").Value = "  -0.17%  "

$ws.Range("D10").Value = "'39.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.21%  "

$ws.Range("D11").Value = "'0.0984"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.78%  "

$ws.Range("D12").Value = "'7.47"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.67%  "

$ws.Range("E13").Value = "  +1.53%  "

$ws.Range("D14").Value = "'2.643.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.95%  "

$ws.Range("D15").Value = "'15.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.22%  "

$ws.Range("D16").Value = "'0.873"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.63%  "

$ws.Range("D17").Value = "'2.298.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.11%  "

$ws.Range("D18").Value = "'42.944.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.48%  "

$ws.Range("E19").Value = "  +3.95%  "

$ws.Range("D20").Value = "'6.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.39%  "

$ws.Range("D21").Value = "'72.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.57%  "

$ws.Range("D22").Value = "'239.20"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.20%  "

$ws.Range("E23").Value = "  +8.58%  "

$ws.Range("D24").Value = "'3.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.53%  "

$ws.Range("D25").Value = "'11.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.00%  "

$ws.Range("E26").Value = "  +0.03%  "

$ws.Range("E27").Value = "  -1.20%  "

$ws.Range("E28").Value = "  -1.32%  "

$ws.Range("D29").Value = "'2.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.37%  "

$ws.Range("D30").Value = "'167.65"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.10%  "

$ws.Range("D31").Value = "'21.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.60%  "

$ws.Range("E32").Value = "  +3.08%  "

$ws.Range("E33").Value = "  +7.30%  "

$ws.Range("E34").Value = "  -1.37%  "

$ws.Range("D35").Value = "'31.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.29%  "

$ws.Range("D36").Value = "'0.127"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.19%  "

$ws.Range("D37").Value = "'4.63"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +11.88%  "

$ws.Range("D38").Value = "'4.83"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.45%  "

$ws.Range("E39").Value = "  -3.11%  "

$ws.Range("E40").Value = "  +11.29%  "

$ws.Range("E41").Value = "  +3.30%  "

$ws.Range("D42").Value = "'5.92"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.64%  "

$ws.Range("D43").Value = "'0.220"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.91%  "

$ws.Range("E44").Value = "  +1.88%  "

$ws.Range("D45").Value = "'62.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.83%  "

$ws.Range("D46").Value = "'4.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.06%  "

$ws.Range("E47").Value = "  +1.69%  "

$ws.Range("D48").Value = "'105.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +11.56%  "

$ws.Range("E49").Value = "  -0.12%  "

$ws.Range("E50").Value = "  +0.04%  "

$ws.Range("E51").Value = "  -0.15%  "
